# Comfenalco Cartagena - "Estado de Cuenta" update
# - Updates Valor Mora total, worker/period counters
# - Replaces the 3 old debt-detail rows with 6 new debt-detail rows
#   (two employees with two periods each, plus the two original rows
#   for the remaining employees), keeping the same look & feel
#   (styles / borders / merged footer) as the original sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header totals (row 11: Valor Mora; row 13: worker/period counts)
# ---------------------------------------------------------------
$ws.Range("E11").Value = 254396

$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 5

# ---------------------------------------------------------------
# 2. Make room for the extra detail rows.
#    Originally rows 16-18 held the 3 detail lines (row 18 carries
#    the "closing" bottom-border style) and rows 23-24 held the
#    signature footer. We need 6 detail rows (16-21), so insert 3
#    blank rows before the old last detail row (18), which pushes
#    that row down to 21 and the footer down to 26-27 automatically,
#    carrying the merged cells with it.
# ---------------------------------------------------------------
$ws.Rows("18:20").Insert()

# Clone the formatting (styles/borders/number formats) of row 17
# (a "normal" detail row) into the 3 new rows.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 3. Fill in the detail rows with the new data.
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo, F=Valor Mora,
#             G=Salario Basico
# ---------------------------------------------------------------

# Row 16: Carmen Cecilia Castro Jimenez - period 2503
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143342065"
$ws.Range("D16").Value = "CARMEN CECILIA CASTRO JIMENEZ"
$ws.Range("E16").Value = "2503"
$ws.Range("F16").Value = 9117
$ws.Range("G16").Value = 6838000

# Row 17: Carmen Cecilia Castro Jimenez - period 2503 (second line)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143342065"
$ws.Range("D17").Value = "CARMEN CECILIA CASTRO JIMENEZ"
$ws.Range("E17").Value = "2503"
$ws.Range("F17").Value = 61196
$ws.Range("G17").Value = 1529924

# Row 18: Luisa Fernanda Villa Julio - period 2502
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1001977563"
$ws.Range("D18").Value = "LUISA FERNANDA VILLA JULIO"
$ws.Range("E18").Value = "2502"
$ws.Range("F18").Value = 61196
$ws.Range("G18").Value = 1529924

# Row 19: Sugey Del Carmen Maturana Rosenstand - period 2206
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "33101883"
$ws.Range("D19").Value = "SUGEY DEL CARMEN MATURANA ROSENSTAND"
$ws.Range("E19").Value = "2206"
$ws.Range("F19").Value = 40133
$ws.Range("G19").Value = 6934400

# Row 20: Joana Marcela Perez - period 1902
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "53124729"
$ws.Range("D20").Value = "JOANA MARCELA PEREZ"
$ws.Range("E20").Value = "1902"
$ws.Range("F20").Value = 41377
$ws.Range("G20").Value = 1075809

# Row 21 (closing row, keeps the bottom-border style): Joana Marcela
# Perez - period 1901
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "53124729"
$ws.Range("D21").Value = "JOANA MARCELA PEREZ"
$ws.Range("E21").Value = "1901"
$ws.Range("F21").Value = 41377
$ws.Range("G21").Value = 1075809
